$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range("D2") "24.632.03"
Set-TextValue $ws.Range("E2") "  -0.39%  "
Set-TextValue $ws.Range("D3") "1.695.66"
Set-TextValue $ws.Range("E3") "  -0.03%  "
Set-TextValue $ws.Range("D4") "1.009"
Set-TextValue $ws.Range("E4") "  +0.91%  "
Set-TextValue $ws.Range("D5") "315.42"
Set-TextValue $ws.Range("E5") "  -0.43%  "
Set-TextValue $ws.Range("D6") "1.008"
Set-TextValue $ws.Range("E6") "  +0.77%  "
Set-TextValue $ws.Range("D7") "0.3926"
Set-TextValue $ws.Range("E7") "  -0.49%  "
Set-TextValue $ws.Range("D8") "0.4057"
Set-TextValue $ws.Range("E8") "  -0.02%  "
Set-TextValue $ws.Range("D9") "1.509"
Set-TextValue $ws.Range("E9") "  +1.39%  "
Set-TextValue $ws.Range("D10") "1.010"
Set-TextValue $ws.Range("E10") "  +0.98%  "
Set-TextValue $ws.Range("D11") "53.13"
Set-TextValue $ws.Range("E11") "  -0.28%  "
Set-TextValue $ws.Range("D12") "0.08768"
Set-TextValue $ws.Range("E12") "  -1.21%  "
Set-TextValue $ws.Range("D13") "7.621"
Set-TextValue $ws.Range("E13") "  +5.39%  "
Set-TextValue $ws.Range("D14") "24.32"
Set-TextValue $ws.Range("E14") "  +2.90%  "
Set-TextValue $ws.Range("D15") "0.00001363"
Set-TextValue $ws.Range("E15") "  +2.99%  "
Set-TextValue $ws.Range("D16") "7.993"
Set-TextValue $ws.Range("E16") "  -0.89%  "
Set-TextValue $ws.Range("D17") "1.695.42"
Set-TextValue $ws.Range("E17") "  +0.12%  "
Set-TextValue $ws.Range("D18") "98.67"
Set-TextValue $ws.Range("E18") "  -1.30%  "
Set-TextValue $ws.Range("D19") "0.07117"
Set-TextValue $ws.Range("E19") "  +1.47%  "
Set-TextValue $ws.Range("D20") "19.80"
Set-TextValue $ws.Range("E20") "  +0.84%  "
Set-TextValue $ws.Range("D21") "7.394"
Set-TextValue $ws.Range("E21") "  +4.71%  "
Set-TextValue $ws.Range("E22") "  +1.01%  "
Set-TextValue $ws.Range("D23") "14.30"
Set-TextValue $ws.Range("E23") "  -0.40%  "
Set-TextValue $ws.Range("D24") "24.622.62"
Set-TextValue $ws.Range("E24") "  -0.35%  "
Set-TextValue $ws.Range("D25") "3.050"
Set-TextValue $ws.Range("E25") "  -6.49%  "
Set-TextValue $ws.Range("D26") "2.357"
Set-TextValue $ws.Range("E26") "  +0.12%  "
Set-TextValue $ws.Range("D27") "22.72"
Set-TextValue $ws.Range("E27") "  -0.05%  "
Set-TextValue $ws.Range("D28") "165.07"
Set-TextValue $ws.Range("E28") "  +1.03%  "
Set-TextValue $ws.Range("D29") "8.469"
Set-TextValue $ws.Range("E29") "  +13.01%  "
Set-TextValue $ws.Range("D30") "137.86"
Set-TextValue $ws.Range("E30") "  +1.18%  "
Set-TextValue $ws.Range("D31") "5.223"
Set-TextValue $ws.Range("E31") "  +0.64%  "
Set-TextValue $ws.Range("D32") "1.884.14"
Set-TextValue $ws.Range("E32") "  +0.16%  "
Set-TextValue $ws.Range("D33") "0.08836"
Set-TextValue $ws.Range("E33") "  +2.71%  "
Set-TextValue $ws.Range("D34") "7.441"
Set-TextValue $ws.Range("E34") "  +3.84%  "
Set-TextValue $ws.Range("D35") "1.050"
Set-TextValue $ws.Range("E35") "  -1.57%  "
Set-TextValue $ws.Range("D36") "1.993"
Set-TextValue $ws.Range("E36") "  +3.16%  "
Set-TextValue $ws.Range("D37") "0.02921"
Set-TextValue $ws.Range("E37") "  +6.61%  "
Set-TextValue $ws.Range("D38") "0.2731"
Set-TextValue $ws.Range("E38") "  -0.70%  "
Set-TextValue $ws.Range("D39") "10.80"
Set-TextValue $ws.Range("E39") "  -6.61%  "
Set-TextValue $ws.Range("D40") "14.30"
Set-TextValue $ws.Range("E40") "  -1.40%  "
Set-TextValue $ws.Range("D41") "0.09144"
Set-TextValue $ws.Range("E41") "  -1.15%  "
Set-TextValue $ws.Range("D42") "0.7865"
Set-TextValue $ws.Range("E42") "  +2.44%  "
Set-TextValue $ws.Range("D43") "1.465"
Set-TextValue $ws.Range("E43") "  -0.52%  "
Set-TextValue $ws.Range("D44") "16.58"
Set-TextValue $ws.Range("E44") "  +3.22%  "
Set-TextValue $ws.Range("D45") "0.7209"
Set-TextValue $ws.Range("E45") "  +0.14%  "
Set-TextValue $ws.Range("D46") "2.565"
Set-TextValue $ws.Range("E46") "  -0.52%  "
Set-TextValue $ws.Range("D47") "4.226"
Set-TextValue $ws.Range("E47") "  +0.17%  "
Set-TextValue $ws.Range("E48") "  +0.74%  "
Set-TextValue $ws.Range("D49") "1.333"
Set-TextValue $ws.Range("E49") "  +0.46%  "
Set-TextValue $ws.Range("D50") "139.31"
Set-TextValue $ws.Range("E50") "  -0.33%  "
Set-TextValue $ws.Range("D51") "91.80"
Set-TextValue $ws.Range("E51") "  +1.61%  "
